# Adds a new "ABNT" column (Z) to the "Classes" worksheet, mirroring the
# existing "ClasseIfc" column (Y): a header in row 1 and the literal text
# "null" in rows 2-8, matching the structure and formatting already used by
# the adjacent CategoriaRvt / ClasseIfc columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classes")

# Clone the formatting of the neighboring "ClasseIfc" column onto the new
# column first (this also extends the sheet's used range to include column
# Z before any values are written, so Excel doesn't re-flow existing row
# heights when the cells below get their content).
$ws.Range("Y1:Y8").Copy()
$ws.Range("Z1:Z8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header for the new column.
$ws.Range("Z1").Value = "ABNT"

# Body rows, same placeholder text as the neighboring columns.
$ws.Range("Z2").Value = "null"
$ws.Range("Z3").Value = "null"
$ws.Range("Z4").Value = "null"
$ws.Range("Z5").Value = "null"
$ws.Range("Z6").Value = "null"
$ws.Range("Z7").Value = "null"
$ws.Range("Z8").Value = "null"

# Keep row heights as originally sized.
$ws.Rows.Item(1).RowHeight = 40.15
$ws.Rows.Item(2).RowHeight = 7.35
$ws.Rows.Item(3).RowHeight = 7.35
$ws.Rows.Item(4).RowHeight = 7.35
$ws.Rows.Item(5).RowHeight = 7.35
$ws.Rows.Item(6).RowHeight = 7.35
$ws.Rows.Item(7).RowHeight = 7.35
$ws.Rows.Item(8).RowHeight = 7.35

# Give the new column a tight width similar to its neighbors (narrow, just
# wide enough for the short placeholder text).
$ws.Columns.Item(26).ColumnWidth = 4.2

# Move the active selection onto the newly added column, as in the source
# workbook (selection moved to Z2 after the edit).
$ws.Range("Z2").Select() | Out-Null
